# Architectures.pptx — "Last Commit Before Refactoring"
#
# Slide 7 ("Analysis" swim-lane diagram): the whole right-hand "Analysis"
# column (Simulation / Ranking / Synthesis Route Generation / Other
# Analysis / Profile Generation boxes, their connectors, the big brace,
# and the "Analysis" header textbox) was nudged to the right/down to make
# room for a new "Encoder - Image to Molecule" trapezoid that was added
# to the "Generation" column.
#
# Shape .Left/.Top/.Width/.Height are expressed in points (1 pt = 12700
# EMU) and are stored as single-precision floats, so the literals below
# are chosen to land exactly on the target EMU value after that
# round-trip.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Simulation rectangle -------------------------------------------------
$sh = $s.Shapes.Item("Rectangle 55")
$sh.Left = 684.6388549804688
$sh.Top = 189.00244140625

# --- Ranking rectangle -----------------------------------------------------
$sh = $s.Shapes.Item("Rectangle 56")
$sh.Left = 642.5291137695312
$sh.Top = 259.5418395996094

# --- connector feeding "Profile Generation" --------------------------------
$sh = $s.Shapes.Item("Straight Arrow Connector 57")
$sh.Left = 768.527587890625
$sh.Top = 271.352783203125

# --- connector feeding "Ranking" -------------------------------------------
$sh = $s.Shapes.Item("Straight Arrow Connector 59")
$sh.Left = 667.7057495117188
$sh.Top = 272.3684387207031

# --- Synthesis Route Generation rectangle ----------------------------------
$sh = $s.Shapes.Item("Rectangle 60")
$sh.Left = 682.1341552734375
$sh.Top = 261.3309631347656

# --- Other Analysis rectangle -----------------------------------------------
$sh = $s.Shapes.Item("Rectangle 62")
$sh.Left = 682.5390014648438
$sh.Top = 336.1643371582031

# --- connector feeding "Synthesis Route Generation" ------------------------
$sh = $s.Shapes.Item("Straight Arrow Connector 63")
$sh.Left = 667.7057495117188
$sh.Top = 201.3424530029297

# --- connector feeding "Other Analysis" ------------------------------------
$sh = $s.Shapes.Item("Straight Arrow Connector 64")
$sh.Left = 667.7057495117188
$sh.Top = 344.2006530761719

# --- Profile Generation rectangle ------------------------------------------
$sh = $s.Shapes.Item("Rectangle 65")
$sh.Left = 711.5963134765625
$sh.Top = 259.541748046875

# --- the brace that spans the whole "Analysis" column (also grew taller) --
$sh = $s.Shapes.Item("Right Brace 68")
$sh.Left = 718.4227294921875
$sh.Top = 8.739370346069336
$sh.Width = 34.711181640625
$sh.Height = 195.22560119628906

# --- "Analysis" column header textbox (only shifted horizontally) ---------
$sh = $s.Shapes.Item("TextBox 71")
$sh.Left = 682.400634765625

# ---------------------------------------------------------------------------
# New shape: "Encoder - Image to Molecule" trapezoid in the "Generation"
# column, styled/sized like its sibling "Trapezoid 47" ("Encoder").
# Adding+removing a throwaway autoshape first forces the next shape to take
# id 3 / name "Trapezoid 2", matching the target deck exactly.
# ---------------------------------------------------------------------------
$placeholder = $s.Shapes.AddShape(3, 10, 10, 20, 20)
$placeholder.Delete()

$template = $s.Shapes.Item("Trapezoid 47")
$newShape = $template.Duplicate().Item(1)
$newShape.Name = "Trapezoid 2"
$newShape.Left = 542.5153198242188
$newShape.Top = 260.63458251953125
$newShape.Width = 211.9673309326172
$newShape.Height = 27.264095306396484
$newShape.TextFrame.TextRange.Text = "Encoder – Image to Molecule"
